$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44243
$ws.Cells.Item(2, 11).Value = 300
$ws.Cells.Item(2, 12).Value = 320
$ws.Cells.Item(2, 13).Value = 310
$ws.Cells.Item(2, 16).Value = 310

# Row 3
$ws.Cells.Item(3, 4).Value = 44243
$ws.Cells.Item(3, 10).Value = 800
$ws.Cells.Item(3, 12).Value = 320
$ws.Cells.Item(3, 13).Value = 310
$ws.Cells.Item(3, 16).Value = 310

# Row 4
$ws.Cells.Item(4, 4).Value = 44202
$ws.Cells.Item(4, 9).Value = "Segunda"
$ws.Cells.Item(4, 10).Value = 1300
$ws.Cells.Item(4, 11).Value = 230
$ws.Cells.Item(4, 12).Value = 250
$ws.Cells.Item(4, 13).Value = 240
$ws.Cells.Item(4, 16).Value = 240

# Row 5
$ws.Cells.Item(5, 4).Value = 44566
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 1000

# Row 6
$ws.Cells.Item(6, 4).Value = 44580
$ws.Cells.Item(6, 10).Value = 1200
$ws.Cells.Item(6, 11).Value = 380
$ws.Cells.Item(6, 12).Value = 400
$ws.Cells.Item(6, 13).Value = 390
$ws.Cells.Item(6, 15).Value = "Región Metropolitana"
$ws.Cells.Item(6, 16).Value = 390

# Row 7
$ws.Cells.Item(7, 4).Value = 44168
$ws.Cells.Item(7, 10).Value = 1700
$ws.Cells.Item(7, 11).Value = 430
$ws.Cells.Item(7, 12).Value = 450
$ws.Cells.Item(7, 13).Value = 440
$ws.Cells.Item(7, 15).Value = "Perú"
$ws.Cells.Item(7, 16).Value = 440

# Row 8
$ws.Cells.Item(8, 4).Value = 44229
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 1200
$ws.Cells.Item(8, 11).Value = 230
$ws.Cells.Item(8, 12).Value = 250
$ws.Cells.Item(8, 13).Value = 240
$ws.Cells.Item(8, 16).Value = 240

# Row 9
$ws.Cells.Item(9, 4).Value = 44589
$ws.Cells.Item(9, 10).Value = 900
$ws.Cells.Item(9, 11).Value = 325
$ws.Cells.Item(9, 12).Value = 350
$ws.Cells.Item(9, 13).Value = 338
$ws.Cells.Item(9, 16).Value = 338

# Row 10
$ws.Cells.Item(10, 4).Value = 44166
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 1700
$ws.Cells.Item(10, 11).Value = 500
$ws.Cells.Item(10, 12).Value = 530
$ws.Cells.Item(10, 13).Value = 515
$ws.Cells.Item(10, 16).Value = 515

# Row 11
$ws.Cells.Item(11, 4).Value = 44253
$ws.Cells.Item(11, 10).Value = 1200
$ws.Cells.Item(11, 11).Value = 270
$ws.Cells.Item(11, 12).Value = 280
$ws.Cells.Item(11, 13).Value = 275
$ws.Cells.Item(11, 15).Value = "Perú"
$ws.Cells.Item(11, 16).Value = 275

# Row 12
$ws.Cells.Item(12, 4).Value = 44176
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 1300
$ws.Cells.Item(12, 11).Value = 350
$ws.Cells.Item(12, 12).Value = 400
$ws.Cells.Item(12, 13).Value = 375
$ws.Cells.Item(12, 16).Value = 375

# Row 13
$ws.Cells.Item(13, 4).Value = 44175
$ws.Cells.Item(13, 11).Value = 400
$ws.Cells.Item(13, 12).Value = 430
$ws.Cells.Item(13, 13).Value = 415
$ws.Cells.Item(13, 16).Value = 415

# Row 14
$ws.Cells.Item(14, 4).Value = 44547
$ws.Cells.Item(14, 10).Value = 1200
$ws.Cells.Item(14, 12).Value = 370
$ws.Cells.Item(14, 13).Value = 360
$ws.Cells.Item(14, 16).Value = 360

# Row 15
$ws.Cells.Item(15, 4).Value = 44217
$ws.Cells.Item(15, 10).Value = 1600
$ws.Cells.Item(15, 12).Value = 350
$ws.Cells.Item(15, 13).Value = 325
$ws.Cells.Item(15, 16).Value = 325

# Row 16
$ws.Cells.Item(16, 4).Value = 44523
$ws.Cells.Item(16, 9).Value = "Segunda"
$ws.Cells.Item(16, 10).Value = 1000
$ws.Cells.Item(16, 11).Value = 550
$ws.Cells.Item(16, 12).Value = 580
$ws.Cells.Item(16, 13).Value = 565
$ws.Cells.Item(16, 16).Value = 565

# Row 17
$ws.Cells.Item(17, 4).Value = 44172
$ws.Cells.Item(17, 9).Value = "Segunda"
$ws.Cells.Item(17, 10).Value = 1600
$ws.Cells.Item(17, 11).Value = 400
$ws.Cells.Item(17, 12).Value = 420
$ws.Cells.Item(17, 13).Value = 410
$ws.Cells.Item(17, 16).Value = 410

# Row 18
$ws.Cells.Item(18, 4).Value = 44160
$ws.Cells.Item(18, 10).Value = 2000
$ws.Cells.Item(18, 11).Value = 500
$ws.Cells.Item(18, 12).Value = 550
$ws.Cells.Item(18, 13).Value = 525
$ws.Cells.Item(18, 16).Value = 525

# Row 19
$ws.Cells.Item(19, 4).Value = 44214
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 1200
$ws.Cells.Item(19, 11).Value = 400
$ws.Cells.Item(19, 12).Value = 450
$ws.Cells.Item(19, 13).Value = 425
$ws.Cells.Item(19, 15).Value = "Perú"
$ws.Cells.Item(19, 16).Value = 425

# Row 21
$ws.Cells.Item(21, 4).Value = 44530
$ws.Cells.Item(21, 10).Value = 1300
$ws.Cells.Item(21, 11).Value = 450
$ws.Cells.Item(21, 12).Value = 480
$ws.Cells.Item(21, 13).Value = 465
$ws.Cells.Item(21, 16).Value = 465

# Row 22
$ws.Cells.Item(22, 4).Value = 44575
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 1200
$ws.Cells.Item(22, 11).Value = 380
$ws.Cells.Item(22, 12).Value = 400
$ws.Cells.Item(22, 13).Value = 390
$ws.Cells.Item(22, 16).Value = 390

# Row 23
$ws.Cells.Item(23, 4).Value = 44162
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 900
$ws.Cells.Item(23, 11).Value = 500
$ws.Cells.Item(23, 12).Value = 550
$ws.Cells.Item(23, 13).Value = 525
$ws.Cells.Item(23, 16).Value = 525

# Row 24
$ws.Cells.Item(24, 4).Value = 44162
$ws.Cells.Item(24, 10).Value = 1200
$ws.Cells.Item(24, 11).Value = 500
$ws.Cells.Item(24, 12).Value = 550
$ws.Cells.Item(24, 13).Value = 525
$ws.Cells.Item(24, 16).Value = 525

# Row 26
$ws.Cells.Item(26, 4).Value = 44603
$ws.Cells.Item(26, 9).Value = "Tercera"
$ws.Cells.Item(26, 10).Value = 300
$ws.Cells.Item(26, 11).Value = 280
$ws.Cells.Item(26, 12).Value = 300
$ws.Cells.Item(26, 13).Value = 290
$ws.Cells.Item(26, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(26, 16).Value = 290

# Row 27
$ws.Cells.Item(27, 4).Value = 44201
$ws.Cells.Item(27, 9).Value = "Segunda"
$ws.Cells.Item(27, 10).Value = 1800
$ws.Cells.Item(27, 11).Value = 250
$ws.Cells.Item(27, 12).Value = 270
$ws.Cells.Item(27, 13).Value = 260
$ws.Cells.Item(27, 16).Value = 260

# Row 28
$ws.Cells.Item(28, 4).Value = 44301
$ws.Cells.Item(28, 10).Value = 900
$ws.Cells.Item(28, 11).Value = 280
$ws.Cells.Item(28, 12).Value = 300
$ws.Cells.Item(28, 13).Value = 290
$ws.Cells.Item(28, 16).Value = 290

# Row 29
$ws.Cells.Item(29, 4).Value = 44231
$ws.Cells.Item(29, 9).Value = "Segunda"
$ws.Cells.Item(29, 10).Value = 200
$ws.Cells.Item(29, 11).Value = 180
$ws.Cells.Item(29, 12).Value = 200
$ws.Cells.Item(29, 13).Value = 190
$ws.Cells.Item(29, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(29, 16).Value = 190

# Row 30
$ws.Cells.Item(30, 4).Value = 44224
$ws.Cells.Item(30, 10).Value = 1200
$ws.Cells.Item(30, 11).Value = 230
$ws.Cells.Item(30, 12).Value = 250
$ws.Cells.Item(30, 13).Value = 240
$ws.Cells.Item(30, 15).Value = "Perú"
$ws.Cells.Item(30, 16).Value = 240

# Row 31
$ws.Cells.Item(31, 4).Value = 44224
$ws.Cells.Item(31, 9).Value = "Segunda"
$ws.Cells.Item(31, 10).Value = 200
$ws.Cells.Item(31, 11).Value = 200
$ws.Cells.Item(31, 12).Value = 230
$ws.Cells.Item(31, 13).Value = 215
$ws.Cells.Item(31, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(31, 16).Value = 215

# Row 32
$ws.Cells.Item(32, 4).Value = 44650
$ws.Cells.Item(32, 10).Value = 1000
$ws.Cells.Item(32, 11).Value = 325
$ws.Cells.Item(32, 12).Value = 350
$ws.Cells.Item(32, 13).Value = 338
$ws.Cells.Item(32, 16).Value = 338

# Row 33
$ws.Cells.Item(33, 4).Value = 44602
$ws.Cells.Item(33, 10).Value = 1300
$ws.Cells.Item(33, 11).Value = 350
$ws.Cells.Item(33, 12).Value = 380
$ws.Cells.Item(33, 13).Value = 365
$ws.Cells.Item(33, 16).Value = 365

# Row 34
$ws.Cells.Item(34, 4).Value = 44602
$ws.Cells.Item(34, 10).Value = 900
$ws.Cells.Item(34, 11).Value = 300
$ws.Cells.Item(34, 12).Value = 330
$ws.Cells.Item(34, 13).Value = 315
$ws.Cells.Item(34, 16).Value = 315

# Row 35
$ws.Cells.Item(35, 4).Value = 44609
$ws.Cells.Item(35, 10).Value = 1200
$ws.Cells.Item(35, 11).Value = 280
$ws.Cells.Item(35, 12).Value = 300
$ws.Cells.Item(35, 13).Value = 290
$ws.Cells.Item(35, 16).Value = 290
